$wb = $excel.ActiveWorkbook

# Locate the sheet currently named "strategy_id-5008"
$wsOld = $wb.Worksheets.Item("strategy_id-5008")

# Create a copy of it, placed immediately after it, before renaming anything
# (this new copy will become "strategy_id-5009")
$wsOld.Copy([System.Reflection.Missing]::Value, $wsOld)

# Rename the original sheet to "strategy_id-5007"
$wsOld.Name = "strategy_id-5007"

# Rename the newly created copy (directly after wsOld) to "strategy_id-5009"
$wsNew = $wb.Worksheets.Item($wsOld.Index + 1)
$wsNew.Name = "strategy_id-5009"
